# Update "想去人数" (F column) values in the "展览" and "全部类型" sheets
# to reflect the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 10022
$wsExhibit.Range("F3").Value = 222
$wsExhibit.Range("F4").Value = 46
$wsExhibit.Range("F5").Value = 606
$wsExhibit.Range("F6").Value = 477

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 10022
$wsAll.Range("F3").Value = 222
$wsAll.Range("F4").Value = 46
$wsAll.Range("F5").Value = 606
$wsAll.Range("F7").Value = 477
